$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.990.44'
$ws.Range("E2").Value = '  +3.85%  '
$ws.Range("D3").Value = '3.351.50'
$ws.Range("E3").Value = '  +9.10%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +7.09%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '621.69'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.86%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +8.41%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.384'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +2.09%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.346.24'
$ws.Range("E10").Value = '  +9.06%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.792'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").Value = '97.776.12'
$ws.Range("E13").Value = '  +3.96%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '35.74'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +6.05%  '
$ws.Range("D15").Value = '3.982.55'
$ws.Range("E15").Value = '  +9.33%  '
$ws.Range("E16").Value = '  +2.12%  '
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").Value = '3.348.86'
$ws.Range("E18").Value = '  +9.37%  '
$ws.Range("E19").Value = '  +1.50%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.74'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.86%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '479.53'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +9.62%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.85'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.66%  '
$ws.Range("E23").Value = '  +9.54%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.10'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +3.51%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '5.67'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.94%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '87.55'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +3.48%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '11.93'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("E29").Value = '  -0.21%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.188'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +6.13%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.247'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("E33").Value = '  +1.30%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.92%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '27.20'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +7.38%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '518.50'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +8.16%  '
$ws.Range("E37").Value = '  -1.62%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '7.28'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -5.13%  '
$ws.Range("E39").Value = '  +2.90%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '24.81'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.16%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.447'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.83%  '
$ws.Range("E42").Value = '  -0.11%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.61'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.72%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.789'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +17.85%  '
$ws.Range("E45").Value = '  +3.89%  '
$ws.Range("E46").Value = '  +0.01%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '160.66'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  +5.71%  '
$ws.Range("E49").Value = '  +6.97%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '45.48'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +4.23%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '4.48'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +6.05%  '
